$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "V2"  = 55.66
    "V3"  = 54.9
    "V4"  = 53.66
    "V5"  = 66.37
    "V6"  = 68.06
    "V7"  = 59.8
    "V8"  = 47.86
    "V9"  = 47.86
    "V11" = 71.95999999999999
    "V15" = 61.56
    "V16" = 56.05
    "V17" = 63.56
    "V18" = 64.5
    "V19" = 61.16
    "V20" = 62.56
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
